$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- nChildren count: 8 -> 6 ---
$ws.Range("B4").Value = "'6"

# --- row 6 (child 0) ---
$ws.Range("B6").Value = "'2"
$ws.Range("C6").Value = "Elwanda  "
$ws.Range("D6").Value = "Cassy  "
$ws.Range("E6").Value = "7.33,-4.6"
$ws.Range("F6").Value = "Tamisha(mother): 0550693864"
$ws.Range("H6").Value = "'28.0"

# --- row 7 (child 1) ---
$ws.Range("B7").Value = "'4"
$ws.Range("C7").Value = "Francisca  "
$ws.Range("D7").Value = "Stevie  "
$ws.Range("E7").Value = "6.22,-0.63"
$ws.Range("F7").Value = "Bernardine(mother): 0561339273"
$ws.Range("G7").Value = "7:06:00"
$ws.Range("H7").Value = "'22.0"

# --- row 8 (child 2) ---
$ws.Range("E8").Value = "9.89,-0.94"
$ws.Range("G8").Value = "7:11:00"
$ws.Range("H8").Value = "'17.0"

# --- row 9 (child 3) ---
$ws.Range("B9").Value = "'12"
$ws.Range("C9").Value = "Frankie  "
$ws.Range("D9").Value = "Flavia  "
$ws.Range("E9").Value = "8.79,0.28"
$ws.Range("F9").Value = "Cyrus(mother): 0522363358"
$ws.Range("G9").Value = "7:13:00"
$ws.Range("H9").Value = "'15.0"

# --- row 10 (child 4) ---
$ws.Range("B10").Value = "'14"
$ws.Range("C10").Value = "Lorinda  "
$ws.Range("D10").Value = "Tyron  "
$ws.Range("E10").Value = "5.04,3.77"
$ws.Range("F10").Value = "Teresa(grandmother): 0558587699"
$ws.Range("G10").Value = "7:19:00"
$ws.Range("H10").Value = "'9.0"

# --- row 11 (child 5) ---
$ws.Range("B11").Value = "'6"
$ws.Range("C11").Value = "Ema  "
$ws.Range("D11").Value = "Ardell  "
$ws.Range("E11").Value = "3.13,1.94"
$ws.Range("F11").Value = "Carley(grandmother): 0533587167"
$ws.Range("G11").Value = "7:23:00"
$ws.Range("H11").Value = "'5.0"

# --- row 12 becomes the "school" summary row (was child 6) ---
$ws.Range("A12").Value = "school"
$ws.Range("B12").Value = "'3"
$ws.Range("C12").Value = "Ironiah"
$ws.Range("D12").Value = "mySchool"
$ws.Range("E12").Value = "0,0"
$ws.Range("F12").Value = "Shir(secretary): 0523345098"
$ws.Range("H12").ClearContents()

# --- row 13 becomes the "cost" summary row (was child 7) ---
$ws.Range("A13").Value = "cost"
$ws.Range("B13").Value = "'25"
$ws.Range("C13:H13").ClearContents()

# --- row 14 becomes the "time" summary row (was the old "school" row) ---
$ws.Range("A14").Value = "time"
$ws.Range("B14").Value = "'28.0"
$ws.Range("C14:G14").ClearContents()

# --- old rows 15 and 16 (cost/time) no longer exist ---
$ws.Range("A15:B16").ClearContents()
